$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new Wins / Losses / Ties columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (bold, bordered, centered/top)
$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Fill in the team record for every data row (2-59)
$lastRow = 59
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 59
    $ws.Cells.Item($r, 31).Value = 103
    $ws.Cells.Item($r, 32).Value = 0
}
